# Commit: "add number 4 in red"
# Add the value 4 to cell D1 and format it with a red font color.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D1")
$cell.Value = 4
$cell.Font.Color = 255   # RGB(255,0,0) -> red, encoded as BGR integer 0x0000FF

# Move the active selection, matching the saved workbook view state.
$ws.Range("G11").Select() | Out-Null
